# Auto-applies the per-cell updates described in the commit diff.
# (GitHub Actions crypto price/volume refresh - Wed Aug 16 03:08:09 UTC 2023)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) and E (Volume(1h)) hold text-like values (e.g. "29.214.29",
# "  -0.44%  ") that must stay literal text, not be coerced to numbers/dates.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '29.214.29'
$ws.Range("E2").Value = '  -0.44%  '
$ws.Range("D3").Value = '1.829.28'
$ws.Range("E3").Value = '  -0.66%  '
$ws.Range("D4").Value = '0.9992'
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").Value = '236.70'
$ws.Range("E5").Value = '  -1.29%  '
$ws.Range("D6").Value = '0.6081'
$ws.Range("E6").Value = '  -3.52%  '
$ws.Range("D8").Value = '0.07106'
$ws.Range("E8").Value = '  -4.67%  '
$ws.Range("D9").Value = '0.2816'
$ws.Range("E9").Value = '  -2.74%  '
$ws.Range("D10").Value = '23.86'
$ws.Range("E10").Value = '  -4.09%  '
$ws.Range("D11").Value = '0.07668'
$ws.Range("E11").Value = '  -0.74%  '
$ws.Range("D12").Value = '1.820.37'
$ws.Range("E12").Value = '  -1.12%  '
$ws.Range("D13").Value = '4.829'
$ws.Range("D14").Value = '0.00001009'
$ws.Range("E14").Value = '  -2.33%  '
$ws.Range("D15").Value = '0.6373'
$ws.Range("E15").Value = '  -5.82%  '
$ws.Range("D16").Value = '2.074.47'
$ws.Range("E16").Value = '  -0.14%  '
$ws.Range("D17").Value = '79.39'
$ws.Range("E17").Value = '  -2.96%  '
$ws.Range("D18").Value = '5.911'
$ws.Range("E18").Value = '  -5.07%  '
$ws.Range("D19").Value = '29.170.90'
$ws.Range("E19").Value = '  -0.51%  '
$ws.Range("D20").Value = '228.33'
$ws.Range("E20").Value = '  -0.26%  '
$ws.Range("E21").Value = '  -4.01%  '
$ws.Range("D22").Value = '1.0000'
$ws.Range("E22").Value = '  +0.05%  '
$ws.Range("D23").Value = '7.037'
$ws.Range("E23").Value = '  -4.55%  '
$ws.Range("D24").Value = '1.001'
$ws.Range("E24").Value = '  +0.18%  '
$ws.Range("D25").Value = '154.52'
$ws.Range("E25").Value = '  -2.22%  '
$ws.Range("D26").Value = '8.082'
$ws.Range("E26").Value = '  -5.12%  '
$ws.Range("D27").Value = '0.1295'
$ws.Range("E27").Value = '  -4.05%  '
$ws.Range("D28").Value = '16.62'
$ws.Range("E28").Value = '  -4.69%  '
$ws.Range("D29").Value = '1.489'
$ws.Range("E29").Value = '  +2.28%  '
$ws.Range("D30").Value = '0.06507'
$ws.Range("E30").Value = '  -6.03%  '
$ws.Range("D31").Value = '1.456'
$ws.Range("E31").Value = '  -2.25%  '
$ws.Range("D32").Value = '3.827'
$ws.Range("E32").Value = '  -5.67%  '
$ws.Range("D33").Value = '3.816'
$ws.Range("E33").Value = '  -6.23%  '
$ws.Range("E34").Value = '  -0.83%  '
$ws.Range("D35").Value = '1.748'
$ws.Range("E35").Value = '  -4.54%  '
$ws.Range("D36").Value = '0.6513'
$ws.Range("E36").Value = '  -7.03%  '
$ws.Range("D37").Value = '2.554'
$ws.Range("E37").Value = '  -1.18%  '
$ws.Range("D38").Value = '2.758'
$ws.Range("E38").Value = '  -2.09%  '
$ws.Range("D39").Value = '1.214.56'
$ws.Range("E39").Value = '  -1.77%  '
$ws.Range("D40").Value = '0.01751'
$ws.Range("E40").Value = '  -5.15%  '
$ws.Range("D41").Value = '6.498'
$ws.Range("E41").Value = '  -4.08%  '
$ws.Range("D42").Value = '0.9325'
$ws.Range("E42").Value = '  -1.18%  '
$ws.Range("D43").Value = '0.9996'
$ws.Range("E43").Value = '  +0.06%  '
$ws.Range("D44").Value = '101.04'
$ws.Range("E44").Value = '  +0.03%  '
$ws.Range("D45").Value = '1.977.78'
$ws.Range("E45").Value = '  -0.28%  '
$ws.Range("D46").Value = '63.33'
$ws.Range("E46").Value = '  -3.09%  '
$ws.Range("D47").Value = '0.00000000119'
$ws.Range("E47").Value = '  -0.53%  '
$ws.Range("D48").Value = '8.586'
$ws.Range("E48").Value = '  -4.55%  '
$ws.Range("E49").Value = '  -5.75%  '
$ws.Range("D50").Value = '0.1079'
$ws.Range("E50").Value = '  -5.66%  '
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").Value = '0.05534'
$ws.Range("E51").Value = '  -2.51%  '
